# Weekly price update: a new "Nectarín" record for
# Agrícola del Norte S.A. de Arica is inserted as row 26, pushing the
# existing rows 26-35 down to 27-36 (their contents are unchanged, only
# their row position shifts). The sheet's used range grows from
# A1:T35 to A1:T36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 26; Excel shifts rows 26:35 down
# to 27:36 and copies formatting (e.g. the date style on column D) from the
# row being pushed down, matching the original workbook's style layout.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44524
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100103
$ws.Range("H26").Value = "Frutos de hueso (carozo)"
$ws.Range("I26").Value = 100103006
$ws.Range("J26").Value = "Nectarín"
$ws.Range("K26").Value = "Early Glo"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 250
$ws.Range("N26").Value = 24000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 24500
$ws.Range("Q26").Value = "$/caja 18 kilos granel"
$ws.Range("R26").Value = "Región de Coquimbo"
$ws.Range("S26").Value = 1361
$ws.Range("T26").Value = 18
